$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $lastPara.Range.Start
$r = $d.Range($insertPoint, $insertPoint)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
      <w:pPr>
        <w:rPr>
          <w:b />
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b />
        </w:rPr>
        <w:t xml:space="preserve">Lesson </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b />
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:contextualSpacing />
      </w:pPr>
      <w:r>
        <w:t>Write a program that accepts one or more numbers as command-line arguments and prints the sum of those numbers to the console (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:r>
        <w:t>stdout</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:r>
        <w:t>).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:contextualSpacing />
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:contextualSpacing />
        <w:rPr>
          <w:b />
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b />
        </w:rPr>
        <w:t>Solution</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:contextualSpacing />
      </w:pPr>
      <w:r>
        <w:tab />
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>var</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t xml:space="preserve"> data = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:r>
        <w:t>process.argv</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:proofErr w:type="spellStart" />
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>var</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t xml:space="preserve"> info = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:r>
        <w:t>data.slice</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:r>
        <w:t>(2);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>function</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t xml:space="preserve"> sum(array){</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>var</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t xml:space="preserve"> result = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:r>
        <w:t>array.reduce</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:r>
        <w:t>(function(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:r>
        <w:t>a,b</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:r>
        <w:t>){</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        a = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>parseInt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t>a);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        b = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>parseInt</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t>b);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        </w:t>
      </w:r>
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>return</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart" />
      <w:r>
        <w:t>a+b</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd" />
      <w:r>
        <w:t>;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    }</w:t>
      </w:r>
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>,0</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t>);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>console.log(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t>result);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:r>
        <w:t>};</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" />
        <w:contextualSpacing />
      </w:pPr>
      <w:proofErr w:type="gramStart" />
      <w:r>
        <w:t>sum(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd" />
      <w:r>
        <w:t>info);</w:t>
      </w:r>
    </w:p>
    
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@

$r.InsertXML($xml)

Write-Host "Paragraphs count after:" $d.Paragraphs.Count
